$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 214
$ws.Range("I12").Value = 217.8
$ws.Range("J12").Value = 195
$ws.Range("K12").Value = 217.8
$ws.Range("L12").Value = 195
$ws.Range("M12").Value = -47.80000000000001
$ws.Range("N12").Value = -535
$ws.Range("H33").Value = 401.33334
$ws.Range("I33").Value = 391.4
$ws.Range("J33").Value = 408.42856
$ws.Range("K33").Value = 391.4
$ws.Range("L33").Value = 408.42856
$ws.Range("M33").Value = -162.4
$ws.Range("N33").Value = -866.4285600000001
$ws.Range("H64").Value = 24994.908
$ws.Range("I64").Value = 28327.334
$ws.Range("J64").Value = 9999
$ws.Range("K64").Value = 28327.334
$ws.Range("L64").Value = 9999
$ws.Range("M64").Value = -28079.334
$ws.Range("N64").Value = -10495
$ws.Range("H67").Value = 24994.908
$ws.Range("I67").Value = 28327.334
$ws.Range("J67").Value = 9999
$ws.Range("K67").Value = 28327.334
$ws.Range("L67").Value = 9999
$ws.Range("M67").Value = -27469.334
$ws.Range("N67").Value = -11715
$ws.Range("H68").Value = 59998.5
$ws.Range("J68").Value = 59998.5
$ws.Range("L68").Value = 59998.5
$ws.Range("N68").Value = -61496.5
$ws.Range("H71").Value = 59998.5
$ws.Range("J71").Value = 59998.5
$ws.Range("L71").Value = 179995.5
$ws.Range("N71").Value = -187483.5
$ws.Range("H116").Value = 1014545.44
$ws.Range("I116").Value = 1192688.9
$ws.Range("J116").Value = 5066
$ws.Range("K116").Value = 1192688.9
$ws.Range("L116").Value = 5066
$ws.Range("M116").Value = -1189246.9
$ws.Range("N116").Value = -11950
$ws.Range("H132").Value = 18098.953
$ws.Range("I132").Value = 19877.895
$ws.Range("K132").Value = 59633.685
$ws.Range("M132").Value = -57103.685
$ws.Range("H135").Value = 2340.8948
$ws.Range("I135").Value = 2340.8948
$ws.Range("K135").Value = 21068.0532
$ws.Range("M135").Value = -18533.0532
$ws.Range("H137").Value = 200668
$ws.Range("I137").Value = 151001
$ws.Range("J137").Value = 300002
$ws.Range("K137").Value = 453003
$ws.Range("L137").Value = 900006
$ws.Range("M137").Value = -450453
$ws.Range("N137").Value = -905106
$ws.Range("H138").Value = 19771.896
$ws.Range("I138").Value = 1700.3704
$ws.Range("J138").Value = 35511.613
$ws.Range("K138").Value = 5101.1112
$ws.Range("L138").Value = 106534.839
$ws.Range("M138").Value = 38.88879999999972
$ws.Range("N138").Value = -116814.839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 48499.5
$ws.Range("J124").Value = 48499.5
$ws.Range("L124").Value = 48499.5
$ws.Range("N124").Value = -58319.5
$ws.Range("H129").Value = 89853.336
$ws.Range("J129").Value = 89853.336
$ws.Range("L129").Value = 89853.336
$ws.Range("N129").Value = -99853.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33334000
$ws.Range("I31").Value = 33334000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 33334000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -33333705
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 33334000
$ws.Range("I34").Value = 33334000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 33334000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -33333798
$ws.Range("N34").ClearContents()
$ws.Range("H63").Value = 39000
$ws.Range("J63").Value = 39000
$ws.Range("L63").Value = 39000
$ws.Range("N63").Value = -40372
$ws.Range("H66").Value = 39000
$ws.Range("J66").Value = 39000
$ws.Range("L66").Value = 117000
$ws.Range("N66").Value = -123864
$ws.Range("H87").Value = 34900
$ws.Range("J87").Value = 34900
$ws.Range("L87").Value = 34900
$ws.Range("N87").Value = -37272
$ws.Range("H90").Value = 34900
$ws.Range("J90").Value = 34900
$ws.Range("L90").Value = 104700
$ws.Range("N90").Value = -116556
$ws.Range("H115").Value = 29642.785
$ws.Range("J115").Value = 29642.785
$ws.Range("L115").Value = 29642.785
$ws.Range("N115").Value = -31992.785
$ws.Range("H134").Value = 2381.6667
$ws.Range("I134").Value = 1640.5264
$ws.Range("K134").Value = 4921.5792
$ws.Range("M134").Value = -2386.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 8499.695
$ws.Range("I32").Value = 7527.778
$ws.Range("J32").Value = 11998.6
$ws.Range("K32").Value = 22583.334
$ws.Range("L32").Value = 35995.8
$ws.Range("M32").Value = -22300.334
$ws.Range("N32").Value = -36561.8
$ws.Range("H46").Value = 834.9
$ws.Range("J46").Value = 1116.3334
$ws.Range("L46").Value = 3349.0002
$ws.Range("N46").Value = -3531.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 33333.332
$ws.Range("J74").Value = 33333.332
$ws.Range("L74").Value = 33333.332
$ws.Range("N74").Value = -35205.332
$ws.Range("H77").Value = 33333.332
$ws.Range("J77").Value = 33333.332
$ws.Range("L77").Value = 99999.99600000001
$ws.Range("N77").Value = -109359.996
$ws.Range("H104").Value = 45890
$ws.Range("J104").Value = 45890
$ws.Range("L104").Value = 45890
$ws.Range("N104").Value = -52878
$ws.Range("H113").Value = 1384.0714
$ws.Range("I113").Value = 1289.75
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 1289.75
$ws.Range("L113").Value = 1950
$ws.Range("M113").Value = 880.25
$ws.Range("N113").Value = -6290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4232.4546
$ws.Range("I7").Value = 4968
$ws.Range("K7").Value = 4968
$ws.Range("M7").Value = -4856
$ws.Range("H126").Value = 4232.4546
$ws.Range("I126").Value = 4968
$ws.Range("K126").Value = 14904
$ws.Range("M126").Value = -12434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 46282.688
$ws.Range("I132").Value = 49238.2
$ws.Range("K132").Value = 147714.6
$ws.Range("M132").Value = -145184.6
